$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1208.849
$ws.Range("I15").Value = 1208.849
$ws.Range("K15").Value = 3626.547
$ws.Range("M15").Value = -3457.547

$ws.Range("H43").Value = 2209.7
$ws.Range("I43").Value = 2519.8
$ws.Range("K43").Value = 2519.8
$ws.Range("M43").Value = -2450.8

$ws.Range("H76").Value = 5540.591
$ws.Range("I76").Value = 4318.0713
$ws.Range("K76").Value = 4318.0713
$ws.Range("M76").Value = -4003.0713

$ws.Range("H79").Value = 5540.591
$ws.Range("I79").Value = 4318.0713
$ws.Range("K79").Value = 4318.0713
$ws.Range("M79").Value = -3226.0713

$ws.Range("H86").Value = 2537.8
$ws.Range("I86").Value = 793.9231
$ws.Range("J86").Value = 3568.2727
$ws.Range("K86").Value = 793.9231
$ws.Range("L86").Value = 3568.2727
$ws.Range("M86").Value = 329.0769
$ws.Range("N86").Value = -5814.2727

$ws.Range("H89").Value = 2537.8
$ws.Range("I89").Value = 793.9231
$ws.Range("J89").Value = 3568.2727
$ws.Range("K89").Value = 3969.6155
$ws.Range("L89").Value = 17841.3635
$ws.Range("M89").Value = 1646.3845
$ws.Range("N89").Value = -29073.3635

$ws.Range("H95").Value = 62499.2
$ws.Range("J95").Value = 70624
$ws.Range("L95").Value = 70624
$ws.Range("N95").Value = -76116

$ws.Range("H98").Value = 2260.6667
$ws.Range("J98").Value = 1874
$ws.Range("L98").Value = 1874
$ws.Range("N98").Value = -4870

$ws.Range("H122").Value = 2260.6667
$ws.Range("J122").Value = 1874
$ws.Range("L122").Value = 5622
$ws.Range("N122").Value = -10522

$ws.Range("H135").Value = 747.3570999999999
$ws.Range("I135").Value = 653.2083
$ws.Range("J135").Value = 1312.25
$ws.Range("K135").Value = 5878.8747
$ws.Range("L135").Value = 11810.25
$ws.Range("M135").Value = -3343.8747
$ws.Range("N135").Value = -16880.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2122.1428
$ws.Range("I14").Value = 4535.2
$ws.Range("J14").Value = 781.55554
$ws.Range("K14").Value = 4535.2
$ws.Range("L14").Value = 781.55554
$ws.Range("M14").Value = -4360.2
$ws.Range("N14").Value = -1131.55554

$ws.Range("H45").Value = 5320.8
$ws.Range("I45").Value = 7501.5
$ws.Range("K45").Value = 7501.5
$ws.Range("M45").Value = -7124.5

$ws.Range("H61").Value = 5935.162
$ws.Range("I61").Value = 6371.88
$ws.Range("K61").Value = 6371.88
$ws.Range("M61").Value = -6159.88

$ws.Range("H110").Value = 879.7692
$ws.Range("I110").Value = 879.7692
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 879.7692
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1165.2308
$ws.Range("N110").ClearContents()

$ws.Range("H136").Value = 5935.162
$ws.Range("I136").Value = 6371.88
$ws.Range("K136").Value = 19115.64
$ws.Range("M136").Value = -16565.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 39199.6
$ws.Range("J35").Value = 39199.6
$ws.Range("L35").Value = 39199.6
$ws.Range("N35").Value = -39819.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2222.5
$ws.Range("I58").Value = 1614.7646
$ws.Range("K58").Value = 1614.7646
$ws.Range("M58").Value = -1411.7646

$ws.Range("H62").Value = 7265.6665
$ws.Range("J62").Value = 7719
$ws.Range("L62").Value = 7719
$ws.Range("N62").Value = -8967

$ws.Range("H65").Value = 7265.6665
$ws.Range("J65").Value = 7719
$ws.Range("L65").Value = 38595
$ws.Range("N65").Value = -44835

$ws.Range("H86").Value = 337435.66
$ws.Range("I86").Value = 503003.5
$ws.Range("K86").Value = 503003.5
$ws.Range("M86").Value = -501880.5

$ws.Range("H89").Value = 337435.66
$ws.Range("I89").Value = 503003.5
$ws.Range("K89").Value = 2515017.5
$ws.Range("M89").Value = -2509401.5

$ws.Range("H134").Value = 1730.1578
$ws.Range("I134").Value = 1629.1538
$ws.Range("K134").Value = 4887.4614
$ws.Range("M134").Value = -2352.4614

$ws.Range("H136").Value = 2222.5
$ws.Range("I136").Value = 1614.7646
$ws.Range("K136").Value = 4844.293799999999
$ws.Range("M136").Value = -2294.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 116
$ws.Range("I44").Value = 76.75
$ws.Range("J44").Value = 155.25
$ws.Range("K44").Value = 230.25
$ws.Range("L44").Value = 465.75
$ws.Range("M44").Value = 167.75
$ws.Range("N44").Value = -1261.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 516433.8
$ws.Range("I21").Value = 1255861.1
$ws.Range("K21").Value = 1255861.1
$ws.Range("M21").Value = -1255688.1

$ws.Range("H30").Value = 516433.8
$ws.Range("I30").Value = 1255861.1
$ws.Range("K30").Value = 1255861.1
$ws.Range("M30").Value = -1255756.1

$ws.Range("H44").Value = 22499.5
$ws.Range("I44").Value = 9999
$ws.Range("K44").Value = 9999
$ws.Range("M44").Value = -9403

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 53799.8
$ws.Range("I4").Value = 41333
$ws.Range("J4").Value = 72500
$ws.Range("K4").Value = 41333
$ws.Range("L4").Value = 72500
$ws.Range("M4").Value = -41220
$ws.Range("N4").Value = -72726

$ws.Range("H28").Value = 53799.8
$ws.Range("I28").Value = 41333
$ws.Range("J28").Value = 72500
$ws.Range("K28").Value = 41333
$ws.Range("L28").Value = 72500
$ws.Range("M28").Value = -41101
$ws.Range("N28").Value = -72964

$ws.Range("H37").Value = 53799.8
$ws.Range("I37").Value = 41333
$ws.Range("J37").Value = 72500
$ws.Range("K37").Value = 41333
$ws.Range("L37").Value = 72500
$ws.Range("M37").Value = -41226
$ws.Range("N37").Value = -72714

$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H46").Value = 2663.1516
$ws.Range("I46").Value = 1111.2
$ws.Range("J46").Value = 3337.913
$ws.Range("K46").Value = 1111.2
$ws.Range("L46").Value = 3337.913
$ws.Range("M46").Value = -923.2
$ws.Range("N46").Value = -3713.913

$ws.Range("H132").Value = 3968.125
$ws.Range("I132").Value = 3499.9092
$ws.Range("J132").Value = 4998.2
$ws.Range("K132").Value = 10499.7276
$ws.Range("L132").Value = 14994.6
$ws.Range("M132").Value = -7969.7276
$ws.Range("N132").Value = -20054.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 78593
$ws.Range("I75").Value = 77895
$ws.Range("K75").Value = 77895
$ws.Range("M75").Value = -76959

$ws.Range("H78").Value = 78593
$ws.Range("I78").Value = 77895
$ws.Range("K78").Value = 233685
$ws.Range("M78").Value = -229005

$ws.Range("H107").Value = 2602.682
$ws.Range("I107").Value = 1656
$ws.Range("K107").Value = 4968
$ws.Range("M107").Value = -3048

$ws.Range("H126").Value = 5370.636
$ws.Range("I126").Value = 5370.636
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16111.908
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -13641.908
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 8386.263000000001
$ws.Range("I132").Value = 8166.357
$ws.Range("K132").Value = 24499.071
$ws.Range("M132").Value = -21969.071

$ws.Range("H136").Value = 9096176
$ws.Range("I136").Value = 11917858
$ws.Range("K136").Value = 35753574
$ws.Range("M136").Value = -35751024
